$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 'The system shall assign a rumour or non-rumor label to a processed Tweet'
$ws.Range("B5").Value = 'To be able to assign a classification label to a particular Tweet after processing has been applied'
$ws.Range("B18").Value = 'The system shall be configurable to account for the needs of the system'
$ws.Range("B20").Value = 'To be able to control the systems services, and enable certain filters, processing techniques'
$ws.Range("B23").Value = 'Configuration files for each service which are deployed at run time'
$ws.Range("B33").Value = 'The system shall provide details on if it is running'
$ws.Range("B35").Value = 'To be able to check if the services are operational '
$ws.Range("B38").Value = 'Various endpoints and checks to ensure services are functional'
$ws.Range("B63").Value = 'The system shall make use of a two-stage classification process'
$ws.Range("B78").Value = 'The system shall accurately store the data to the database'
$ws.Range("B93").Value = 'The system shall be able to link back a Tweet to a particular user'
$ws.Range("B108").Value = 'The system shall be able to link back a Tweet to a particular hashtag'
$ws.Range("B123").Value = 'The system shall be configurable with the filter list'
$ws.Range("B138").Value = 'The system shall only add Tweets that are not a retweet to the Queue'
$ws.Range("B153").Value = 'The system shall be robust enough to restart on failure'
$ws.Range("B168").Value = 'The system shall ensure that the message is valid from the queue'
$ws.Range("B183").Value = 'The system shall be able to report on successes/failures'
$ws.Range("B198").Value = 'The system shall remove/discard any items that are not a Tweet object'
$ws.Range("B213").Value = 'The system shall only handle Tweet objects from the queue reader'
$ws.Range("B228").Value = 'The system shall only handle Tweet objects from the queue reader'
$ws.Range("B243").Value = 'The system shall be able to receive user requests'
$ws.Range("B258").Value = 'The system shall be able to display status information when user navigates to homepage'
$ws.Range("B273").Value = 'The system shall only display relevant information to the status of the system'
$ws.Range("B288").Value = 'The system shall be free from SQL injection attempts'
$ws.Range("B303").Value = 'The system shall be secure'
$ws.Range("B318").Value = 'The system shall be able to respond to user requests within 15 seconds'
$ws.Range("B333").Value = 'The system shall return responses in JSON format'
$ws.Range("B363").Value = 'The system shall provide help to the users'
$ws.Range("B378").Value = 'The system shall allow users to report particular terms'
$ws.Range("B393").Value = 'The system shall be pleasing to the eye'

$ws.Range("B258:F259").Select()
